$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking values are not
# auto-converted to numbers (and precision-altered) by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '52.247.00'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.828.44'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '359.53'
$ws.Range('E5').Value = '  +3.76%  '
$ws.Range('D6').Value = '112.39'
$ws.Range('E6').Value = '  -3.05%  '
$ws.Range('D7').Value = '0.573'
$ws.Range('E7').Value = '  +4.30%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.601'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').Value = '0.0870'
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '19.92'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '7.78'
$ws.Range('E14').Value = '  -1.12%  '
$ws.Range('D15').Value = '3.272.12'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').Value = '2.833.99'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = '0.929'
$ws.Range('E17').Value = '  +4.23%  '
$ws.Range('D18').Value = '52.134.89'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('E19').Value = '  +3.78%  '
$ws.Range('D20').Value = '3.17'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').Value = '0.0₃0999'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Value = '272.88'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('D24').Value = '70.61'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '2.82'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('D26').Value = '27.03'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '10.35'
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('D29').Value = '2.25'
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').Value = '0.0488'
$ws.Range('E30').Value = '  +6.53%  '
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').Value = '35.13'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').Value = '52.10'
$ws.Range('E33').Value = '  +4.16%  '
$ws.Range('E34').Value = '  +3.55%  '
$ws.Range('D35').Value = '5.63'
$ws.Range('E35').Value = '  +14.11%  '
$ws.Range('D36').Value = '0.0852'
$ws.Range('E36').Value = '  +2.45%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('D39').Value = '2.05'
$ws.Range('E39').Value = '  -2.69%  '
$ws.Range('D40').Value = '18.45'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').Value = '0.118'
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('D42').Value = '127.61'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').Value = '2.55'
$ws.Range('E43').Value = '  -1.92%  '
$ws.Range('D44').Value = '23.19'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '2.29'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '2.094.01'
$ws.Range('E46').Value = '  +1.74%  '
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').Value = '2.29'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').Value = '5.92'
$ws.Range('E49').Value = '  +6.52%  '
$ws.Range('D50').Value = '0.970'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = '9.19'
$ws.Range('E51').Value = '  +2.87%  '

# Restore original (default) cell style so no stray style/number-format
# is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
